{"js": "// Applies the four text replacements described by the diff:\n//   1. \"2024-12-24\" -> \"2024-12-25\"                                                   (Fecha)\n//   2. \"Tecnolog\u00eda Superior en Automatizaci\u00f3n e Instrumentaci\u00f3n\"\n//        -> \"Tecnolog\u00eda Superior en Redes y Telecomunicaciones\"                       (Carrera)\n//   3. \"0992783016\" -> \"0998035014\"                                                   (Celular)\n//   4. \"bhh\" -> \"hjghj\"                                                               (Asunto body text)\n\nconst replacements = [\n  [\"2024-12-24\", \"2024-12-25\"],\n  [\"Tecnolog\u00eda Superior en Automatizaci\u00f3n e Instrumentaci\u00f3n\", \"Tecnolog\u00eda Superior en Redes y Telecomunicaciones\"],\n  [\"0992783016\", \"0998035014\"],\n  [\"bhh\", \"hjghj\"],\n];\n\nfor (const [findText, newText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the four text replacements described by the diff:\n#   1. \"2024-12-24\" -> \"2024-12-25\"                                                   (Fecha)\n#   2. \"Tecnolog\u00eda Superior en Automatizaci\u00f3n e Instrumentaci\u00f3n\"\n#        -> \"Tecnolog\u00eda Superior en Redes y Telecomunicaciones\"                       (Carrera)\n#   3. \"0992783016\" -> \"0998035014\"                                                   (Celular)\n#   4. \"bhh\" -> \"hjghj\"                                                               (Asunto body text)\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @{ Find = \"2024-12-24\"; Replace = \"2024-12-25\" },\n    @{ Find = \"Tecnolog\u00eda Superior en Automatizaci\u00f3n e Instrumentaci\u00f3n\"; Replace = \"Tecnolog\u00eda Superior en Redes y Telecomunicaciones\" },\n    @{ Find = \"0992783016\"; Replace = \"0998035014\" },\n    @{ Find = \"bhh\"; Replace = \"hjghj\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $r.Replace, $wdReplaceAll)\n}\n"}
